$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row to the table (ListObject) so the table range, dimension and
# formatting all extend automatically, mirroring what Excel does when a
# row is typed in directly below an existing Excel Table.
$table = $ws.ListObjects.Item("Table1")
$lastRow = $table.ListRows.Item($table.ListRows.Count)
$newRow = $table.ListRows.Add()

# Match the cell formatting used by the previous data row, the same way
# Excel extends formatting for a newly typed table row.
$lastRow.Range.Copy()
$newRow.Range.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$newRow.Range.EntireRow.RowHeight = $lastRow.Range.EntireRow.RowHeight

$dataRange = $newRow.Range

$dataRange.Item(1, 1).Value = 45907
$dataRange.Item(1, 2).Value = 150
$dataRange.Item(1, 3).Value = 776
$dataRange.Item(1, 4).Value = 0
$dataRange.Item(1, 5).Value = 0
$dataRange.Item(1, 6).Value = 1012
$dataRange.Item(1, 7).Value = "N/A"

# Move selection to reflect where the cursor ends up after entering the
# last value of the new row (one row below, same column as last edit).
$ws.Range("J55").Select()
